$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill rows 15-31 (columns A-H) with the same stable value used throughout the table
$value = 0.0026
for ($r = 15; $r -le 31; $r++) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $value
    }
}

# Update selection / view to match the saved state: active cell A14, selection A14:H31,
# with the window scrolled so row 28 is the top visible row.
$ws.Range("A14:H31").Select()
$wb.Windows.Item(1).ScrollRow = 28
